# Add a new weekly data point for "Macroferia Regional de Talca - Repollo".
# A new row is inserted at row 183 (pushing the former rows 183-191 down to
# 184-192), and the new row 183 is populated with the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 183, shifting existing data down.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44509
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100112006
$ws.Range("G183").Value = "Repollo"
$ws.Range("H183").Value = "Crespo record"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 2000
$ws.Range("K183").Value = 900
$ws.Range("L183").Value = 900
$ws.Range("M183").Value = 900
$ws.Range("N183").Value = "`$/unidad"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 900
$ws.Range("Q183").Value = 1
$ws.Range("R183").Value = "Hortaliza"
